# Capacity Supply Curve.xlsx - edit script
# 1. Recalibrate the "Cost multiplier vs Share of Existing Capacity Built" curve
#    (sheet CSC-CSCCCMvSoECBtY) to the new, rounder values.
# 2. Recalibrate the "Share of Cost Effective Capacity Built in a Single Year"
#    table (sheet CSC-CSCSoCECBiaSY): rows 2-15 and 18-25 (cols B:AE) -> 0.4
#    (rows 16/17, already 0, are left untouched).
# 3. Update sheet selections / active sheet to match the saved UI state.

$wb = $excel.ActiveWorkbook

# --- Sheet: CSC-CSCCCMvSoECBtY ---------------------------------------------
$wsCCM = $wb.Worksheets.Item("CSC-CSCCCMvSoECBtY")
$wsCCM.Activate()

$wsCCM.Range("C2").Value = 1.1000000000000001
$wsCCM.Range("D2").Value = 1.3
$wsCCM.Range("E2").Value = 1.6
$wsCCM.Range("F2").Value = 2
$wsCCM.Range("G2").Value = 2.5
$wsCCM.Range("H2").Value = 3.1
$wsCCM.Range("I2").Value = 3.8
$wsCCM.Range("J2").Value = 4.5999999999999996
$wsCCM.Range("K2").Value = 5.5
$wsCCM.Range("L2").Value = 6.4999999999999991
$wsCCM.Range("M2").Value = 7.5999999999999979
$wsCCM.Range("N2").Value = 8.8000000000000007

$wsCCM.Range("C3:O7").Select()

# --- Sheet: CSC-CSCSoCECBiaSY -----------------------------------------------
$wsSoC = $wb.Worksheets.Item("CSC-CSCSoCECBiaSY")
$wsSoC.Activate()

$wsSoC.Range("B2:AE15").Value = 0.4
$wsSoC.Range("B18:AE25").Value = 0.4

$wsSoC.Range("B18:AE25").Select()

# --- Sheet: About (becomes the active/selected tab) -------------------------
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Activate()
$wsAbout.Range("A4").Select()
